# Update loading_percent results for Case_1_15 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16.66518678398402
$ws.Range("D2").Value = 11.33933034893824
$ws.Range("E2").Value = 17.62885846872599
$ws.Range("F2").Value = 27.62294620771841
$ws.Range("G2").Value = 25.47447013607522
$ws.Range("H2").Value = 13.41207476211121
$ws.Range("I2").Value = 24.83148595105494
$ws.Range("J2").Value = 12.05773791558464
$ws.Range("L2").Value = 10.79510665838266
$ws.Range("M2").Value = 16.19514583551026
$ws.Range("O2").Value = 19.98460378381252

# Row 3
$ws.Range("B3").Value = 16.17437411603451
$ws.Range("D3").Value = 11.37644539660668
$ws.Range("E3").Value = 17.59821573736666
$ws.Range("F3").Value = 27.77336721179054
$ws.Range("G3").Value = 25.49446655214466
$ws.Range("H3").Value = 13.46179234910285
$ws.Range("I3").Value = 24.97839359500082
$ws.Range("J3").Value = 12.02950035231013
$ws.Range("L3").Value = 10.54801856234324
$ws.Range("M3").Value = 15.90648031536592
$ws.Range("O3").Value = 20.05552231726789

# Row 4
$ws.Range("B4").Value = 15.86494060245285
$ws.Range("D4").Value = 11.4007473088732
$ws.Range("E4").Value = 17.58158597191007
$ws.Range("F4").Value = 27.87383025980564
$ws.Range("G4").Value = 25.51783447980701
$ws.Range("H4").Value = 13.49499600592076
$ws.Range("I4").Value = 25.07486455153248
$ws.Range("J4").Value = 12.01323661427007
$ws.Range("L4").Value = 10.39301836073373
$ws.Range("M4").Value = 15.72646210000873
$ws.Range("O4").Value = 20.10461729050217

# Row 5
$ws.Range("B5").Value = 15.73696961883554
$ws.Range("D5").Value = 11.41103193706514
$ws.Range("E5").Value = 17.57536489309905
$ws.Range("F5").Value = 27.91680065858917
$ws.Range("G5").Value = 25.53013494578788
$ws.Range("H5").Value = 13.50919905606407
$ws.Range("I5").Value = 25.1157529297159
$ws.Range("J5").Value = 12.00688354943168
$ws.Range("L5").Value = 10.3290977722059
$ws.Range("M5").Value = 15.65247777174396
$ws.Range("O5").Value = 20.12601470569147

# Row 6
$ws.Range("B6").Value = 15.71561171833947
$ws.Range("D6").Value = 11.412762754402
$ws.Range("E6").Value = 17.57436562677066
$ws.Range("F6").Value = 27.9240583189793
$ws.Range("G6").Value = 25.53234483458089
$ws.Range("H6").Value = 13.51159804197025
$ws.Range("I6").Value = 25.12263756170648
$ws.Range("J6").Value = 12.00584533544964
$ws.Range("L6").Value = 10.31844021453109
$ws.Range("M6").Value = 15.6401571252749
$ws.Range("O6").Value = 20.1296515826426

# Row 7
$ws.Range("B7").Value = 15.8632221154598
$ws.Range("D7").Value = 11.40088446546466
$ws.Range("E7").Value = 17.58149981427646
$ws.Range("F7").Value = 27.87440155979139
$ws.Range("G7").Value = 25.51798913713245
$ws.Range("H7").Value = 13.49518483225471
$ws.Range("I7").Value = 25.07540960752341
$ws.Range("J7").Value = 12.01314981729781
$ws.Range("L7").Value = 10.39215927439491
$ws.Range("M7").Value = 15.72546675613379
$ws.Range("O7").Value = 20.10490023890187

# Row 8
$ws.Range("B8").Value = 16.49771957054795
$ws.Range("D8").Value = 11.35181411855682
$ws.Range("E8").Value = 17.61784319011556
$ws.Range("F8").Value = 27.67312525254807
$ws.Range("G8").Value = 25.47905746717138
$ws.Range("H8").Value = 13.42866131275197
$ws.Range("I8").Value = 24.88083810185284
$ws.Range("J8").Value = 12.04778133424549
$ws.Range("L8").Value = 10.71063364414251
$ws.Range("M8").Value = 16.09623308239618
$ws.Range("O8").Value = 20.0079015346274

# Row 9
$ws.Range("B9").Value = 17.67170637807492
$ws.Range("D9").Value = 11.26755223052
$ws.Range("E9").Value = 17.70616449609591
$ws.Range("F9").Value = 27.34303140222436
$ws.Range("G9").Value = 25.49106017480511
$ws.Range("H9").Value = 13.31948492779786
$ws.Range("I9").Value = 24.54906162372647
$ws.Range("J9").Value = 12.12401865041901
$ws.Range("L9").Value = 11.30620805988894
$ws.Range("M9").Value = 16.7982108740103
$ws.Range("O9").Value = 19.86194142331031

# Row 10
$ws.Range("B10").Value = 18.48395346465076
$ws.Range("D10").Value = 11.21288269939104
$ws.Range("E10").Value = 17.78105286292928
$ws.Range("F10").Value = 27.14031595821703
$ws.Range("G10").Value = 25.55402484776882
$ws.Range("H10").Value = 13.25229230421684
$ws.Range("I10").Value = 24.3357040568023
$ws.Range("J10").Value = 12.1848325357228
$ws.Range("L10").Value = 11.72256834979141
$ws.Range("M10").Value = 17.29475568649626
$ws.Range("O10").Value = 19.78195725524263

# Row 11
$ws.Range("B11").Value = 18.84121095189112
$ws.Range("D11").Value = 11.18957178987219
$ws.Range("E11").Value = 17.81720264631264
$ws.Range("F11").Value = 27.05682744199318
$ws.Range("G11").Value = 25.5944246012327
$ws.Range("H11").Value = 13.22456182754242
$ws.Range("I11").Value = 24.24525675141362
$ws.Range("J11").Value = 12.2134802879432
$ws.Range("L11").Value = 11.90670687099754
$ws.Range("M11").Value = 17.51571481753188
$ws.Range("O11").Value = 19.75154106048748

# Row 12
$ws.Range("B12").Value = 18.97463877664989
$ws.Range("D12").Value = 11.18096775571715
$ws.Range("E12").Value = 17.83118305572185
$ws.Range("F12").Value = 27.02647489980705
$ws.Range("G12").Value = 25.61140842717722
$ws.Range("H12").Value = 13.21446964307589
$ws.Range("I12").Value = 24.21195875207096
$ws.Range("J12").Value = 12.22446444195231
$ws.Range("L12").Value = 11.97562784204798
$ws.Range("M12").Value = 17.59861930517901
$ws.Range("O12").Value = 19.74088543657532

# Row 13
$ws.Range("B13").Value = 18.94598678488108
$ws.Range("D13").Value = 11.18281086991525
$ws.Range("E13").Value = 17.82815930347545
$ws.Range("F13").Value = 27.03295557992372
$ws.Range("G13").Value = 25.60767581920236
$ws.Range("H13").Value = 13.21662498144733
$ws.Range("I13").Value = 24.21908768070349
$ws.Range("J13").Value = 12.22209285392112
$ws.Range("L13").Value = 11.96082119584366
$ws.Range("M13").Value = 17.58079941576025
$ws.Range("O13").Value = 19.74314190357184

# Row 14
$ws.Range("B14").Value = 18.85222589744287
$ws.Range("D14").Value = 11.18885946014785
$ws.Range("E14").Value = 17.81834703343414
$ws.Range("F14").Value = 27.05430496887427
$ws.Range("G14").Value = 25.5957881360616
$ws.Range("H14").Value = 13.22372334016143
$ws.Range("I14").Value = 24.2424982054641
$ws.Range("J14").Value = 12.21438126981667
$ws.Range("L14").Value = 11.91239351790142
$ws.Range("M14").Value = 17.52255114257824
$ws.Range("O14").Value = 19.75064711187716

# Row 15
$ws.Range("B15").Value = 18.79454997127401
$ws.Range("D15").Value = 11.1925934527126
$ws.Range("E15").Value = 17.81237441876082
$ws.Range("F15").Value = 27.06754676043667
$ws.Range("G15").Value = 25.58872585401941
$ws.Range("H15").Value = 13.22812454868831
$ws.Range("I15").Value = 24.25696192124444
$ws.Range("J15").Value = 12.20967523082751
$ws.Range("L15").Value = 11.88262345869171
$ws.Range("M15").Value = 17.48677065834127
$ws.Range("O15").Value = 19.75535668206669

# Row 16
$ws.Range("B16").Value = 18.46035096845531
$ws.Range("D16").Value = 11.2144374140974
$ws.Range("E16").Value = 17.77873160179156
$ws.Range("F16").Value = 27.14594841426118
$ws.Range("G16").Value = 25.55162082462491
$ws.Range("H16").Value = 13.25416170217443
$ws.Range("I16").Value = 24.34174817352372
$ws.Range("J16").Value = 12.18297965393651
$ws.Range("L16").Value = 11.7104240416388
$ws.Range("M16").Value = 17.28021122379683
$ws.Range("O16").Value = 19.78406548580448

# Row 17
$ws.Range("B17").Value = 18.25212316549113
$ws.Range("D17").Value = 11.22823657342457
$ws.Range("E17").Value = 17.7586204821579
$ws.Range("F17").Value = 27.19628618260119
$ws.Range("G17").Value = 25.53186651834469
$ws.Range("H17").Value = 13.27086165760743
$ws.Range("I17").Value = 24.3954561206093
$ws.Range("J17").Value = 12.16685072191037
$ws.Range("L17").Value = 11.60339890767324
$ws.Range("M17").Value = 17.15218919270428
$ws.Range("O17").Value = 19.80320916014979

# Row 18
$ws.Range("B18").Value = 18.13120932787215
$ws.Range("D18").Value = 11.23632022973467
$ws.Range("E18").Value = 17.7472496383475
$ws.Range("F18").Value = 27.22606017649036
$ws.Range("G18").Value = 25.52161148713657
$ws.Range("H18").Value = 13.28073389160713
$ws.Range("I18").Value = 24.4269695656227
$ws.Range("J18").Value = 12.15766666885814
$ws.Range("L18").Value = 11.54134801951551
$ws.Range("M18").Value = 17.07809444102746
$ws.Range("O18").Value = 19.81478170824224

# Row 19
$ws.Range("B19").Value = 18.09007619280171
$ws.Range("D19").Value = 11.2390824453543
$ws.Range("E19").Value = 17.74343366833153
$ws.Range("F19").Value = 27.23628191785981
$ws.Range("G19").Value = 25.51832956931132
$ws.Range("H19").Value = 13.28412226438943
$ws.Range("I19").Value = 24.43774625903788
$ws.Range("J19").Value = 12.15457322489722
$ws.Range("L19").Value = 11.52025562960559
$ws.Range("M19").Value = 17.05293014185183
$ws.Range("O19").Value = 19.81879630458083

# Row 20
$ws.Range("B20").Value = 18.27440879479165
$ws.Range("D20").Value = 11.22675244773913
$ws.Range("E20").Value = 17.76074106273501
$ws.Range("F20").Value = 27.19084261050681
$ws.Range("G20").Value = 25.53385484326788
$ws.Range("H20").Value = 13.26905629263626
$ws.Range("I20").Value = 24.38967442424051
$ws.Range("J20").Value = 12.16855810086758
$ws.Range("L20").Value = 11.61484330830185
$ws.Range("M20").Value = 17.16586539827303
$ws.Range("O20").Value = 19.80111312848937

# Row 21
$ws.Range("B21").Value = 18.87981689337663
$ws.Range("D21").Value = 11.18707678869266
$ws.Range("E21").Value = 17.82122129667763
$ws.Range("F21").Value = 27.04799980418712
$ws.Range("G21").Value = 25.59923415553394
$ws.Range("H21").Value = 13.22162727935977
$ws.Range("I21").Value = 24.23559610234184
$ws.Range("J21").Value = 12.2166427059825
$ws.Range("L21").Value = 11.92664022404021
$ws.Range("M21").Value = 17.53968137632786
$ws.Range("O21").Value = 19.74841921677335

# Row 22
$ws.Range("B22").Value = 19.26461898893097
$ws.Range("D22").Value = 11.16244765065317
$ws.Range("E22").Value = 17.86244239619678
$ws.Range("F22").Value = 26.96200754841799
$ws.Range("G22").Value = 25.65178196184081
$ws.Range("H22").Value = 13.19301258647804
$ws.Range("I22").Value = 24.14044948604561
$ws.Range("J22").Value = 12.24885810372903
$ws.Range("L22").Value = 12.12568540542606
$ws.Range("M22").Value = 17.77949016963447
$ws.Range("O22").Value = 19.7190086209215

# Row 23
$ws.Range("B23").Value = 19.06026692427655
$ws.Range("D23").Value = 11.17547389241309
$ws.Range("E23").Value = 17.84028971618139
$ws.Range("F23").Value = 27.00722692162454
$ws.Range("G23").Value = 25.62284036113407
$ws.Range("H23").Value = 13.20806643854978
$ws.Range("I23").Value = 24.1907223132526
$ws.Range("J23").Value = 12.23159376681515
$ws.Range("L23").Value = 12.0199000424836
$ws.Range("M23").Value = 17.65193057055576
$ws.Range("O23").Value = 19.73424432377765

# Row 24
$ws.Range("B24").Value = 18.26433719475778
$ws.Range("D24").Value = 11.22742295262111
$ws.Range("E24").Value = 17.75978175172737
$ws.Range("F24").Value = 27.19330105173251
$ws.Range("G24").Value = 25.53295248861438
$ws.Range("H24").Value = 13.26987165340842
$ws.Range("I24").Value = 24.39228634810976
$ws.Range("J24").Value = 12.1677859181736
$ws.Range("L24").Value = 11.60967091277116
$ws.Range("M24").Value = 17.15968391726769
$ws.Range("O24").Value = 19.80205897979304

# Row 25
$ws.Range("B25").Value = 17.36244963545883
$ws.Range("D25").Value = 11.289072241054
$ws.Range("E25").Value = 17.68048553801768
$ws.Range("F25").Value = 27.42537051468151
$ws.Range("G25").Value = 25.47831346616054
$ws.Range("H25").Value = 13.34673717703155
$ws.Range("I25").Value = 24.63348190897597
$ws.Range("J25").Value = 12.10253077902979
$ws.Range("L25").Value = 11.14857249525924
$ws.Range("M25").Value = 16.61140840084509
$ws.Range("O25").Value = 19.89666038451794

